$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = 350
$ws.Range("B42").Value = 0.33

$ws.Range("B43").Select()
